$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Event Data")
$ws.Range("C95").Value = "SKIP"
Write-Host "C95 after set:" $ws.Range("C95").Text
